# "agregamos color a gato"
# Insert a new "String color2;" declaration line right after the
# "int numeroVidas;" line, moving the _GoBack bookmark into the new
# paragraph so it sits right after "String color2" and before the ";"
# (i.e. in the same relative spot it occupied before: right before the
# trailing semicolon of the line it was attached to).

$d = $word.ActiveDocument

# Find the paragraph whose text is exactly "int numeroVidas;"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "int numeroVidas;") {
        $target = $p
    }
}

# Insert a new, empty paragraph right after it and fill it with the new
# field declaration.
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()
$newText = "String color2;"
$newPara.Range.Text = $newText

# The _GoBack bookmark originally sat collapsed at the end of
# "int numeroVidas;" (right before that paragraph's mark). Relocate it so
# it sits collapsed right after "String color2" and before the ";" in the
# freshly-inserted paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$prefixLen = "String color2".Length
$bmPos = $newPara.Range.Start + $prefixLen
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
